$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Grade (F3) changes from "L3" to "L4"; Technology (G3) changes from "JavaScript" to "Java"
$ws.Range("F3").Value = "L4"
$ws.Range("G3").Value = "Java"

# Update selection to G3
$ws.Range("G3").Select()
